# Updates crypto price/volume data to match the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.740.29"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "'3.497.63"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'602.25"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("D6").Value = "'147.31"
$ws.Range("E6").Value = "  -3.28%  "
$ws.Range("D7").Value = "'3.495.49"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.66%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'7.79"
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -1.34%  "
$ws.Range("D14").Value = "'4.088.20"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'31.30"
$ws.Range("E15").Value = "  -4.54%  "
$ws.Range("D16").Value = "'3.503.47"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'66.745.16"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "'10.56"
$ws.Range("E19").Value = "  +7.15%  "
$ws.Range("D20").Value = "'6.37"
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").Value = "'15.34"
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D22").Value = "'433.75"
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "'0.609"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("D24").Value = "'79.78"
$ws.Range("E24").Value = "  +2.02%  "
$ws.Range("D25").Value = "'3.634.67"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("E28").Value = "  -5.97%  "
$ws.Range("D29").Value = "'9.78"
$ws.Range("E29").Value = "  -3.24%  "
$ws.Range("D30").Value = "'8.24"
$ws.Range("E30").Value = "  -6.96%  "
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("E32").Value = "  -4.31%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("D35").Value = "'25.30"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("D36").Value = "'3.490.24"
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("D38").Value = "'5.88"
$ws.Range("E38").Value = "  -5.35%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("D42").Value = "'0.0890"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "'169.70"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("E44").Value = "  -9.50%  "
$ws.Range("D45").Value = "'5.41"
$ws.Range("E45").Value = "  -1.52%  "
$ws.Range("D46").Value = "'0.897"
$ws.Range("E46").Value = "  +1.70%  "
$ws.Range("D47").Value = "'29.10"
$ws.Range("E47").Value = "  -4.07%  "
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("D49").Value = "'1.32"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'7.46"
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("E51").Value = "  -4.15%  "
